$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three day-blocks (19-03 Status/Time and 20-03 old position) that
# are no longer needed, collapsing the sheet from columns A:I down to A:E.
# Columns F:I (19-03-2025 Status/Time, 20-03-2025 Status/Time) get removed;
# the surviving "20-03-2025" columns end up as D/E.
$ws.Range("F1:I6").Delete()

# Header row: rename the remaining date columns to 20-03-2025
$ws.Range("D1").Value = "20-03-2025 Status"
$ws.Range("E1").Value = "20-03-2025 Time"

# Row 2 (student 4201)
$ws.Range("B2").Value = "muhammad"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "P"
$ws.Range("E2").Value = "13:48:30"

# Row 3 (student 4202)
$ws.Range("B3").Value = "Ayat Eman"
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = "13:48:35"

# Row 4 (student 4203)
$ws.Range("B4").Value = "zaira Malik"
$ws.Range("E4").Value = "13:48:26"

# Row 5 (student 4204)
$ws.Range("B5").Value = "sukeena Ali"
$ws.Range("C5").Value = 1
$ws.Range("E5").Value = "13:48:21"

# Row 6 (student 4205)
$ws.Range("B6").Value = "Ameer Abbas"
$ws.Range("C6").ClearContents()
